$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 607.1739
$ws.Range("J17").Value = 607.1739
$ws.Range("L17").Value = 1821.5217
$ws.Range("N17").Value = -2157.5217

$ws.Range("H40").Value = 1471.5555
$ws.Range("I40").Value = 1224.8334
$ws.Range("K40").Value = 1224.8334
$ws.Range("M40").Value = -1049.8334

$ws.Range("H80").Value = 600.9286
$ws.Range("I80").Value = 546
$ws.Range("J80").Value = 610.0833
$ws.Range("K80").Value = 1638
$ws.Range("L80").Value = 1830.2499
$ws.Range("M80").Value = -640
$ws.Range("N80").Value = -3826.2499

$ws.Range("H83").Value = 600.9286
$ws.Range("I83").Value = 546
$ws.Range("J83").Value = 610.0833
$ws.Range("K83").Value = 4914
$ws.Range("L83").Value = 5490.7497
$ws.Range("M83").Value = 78
$ws.Range("N83").Value = -15474.7497

$ws.Range("H100").Value = 2199
$ws.Range("I100").Value = 2248.75
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2248.75
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1707.75
$ws.Range("N100").Value = -3082

$ws.Range("H121").Value = 1682
$ws.Range("J121").Value = 1727.5
$ws.Range("L121").Value = 5182.5
$ws.Range("N121").Value = -8676.5

$ws.Range("H129").Value = 905.2098999999999
$ws.Range("J129").Value = 960.58105
$ws.Range("L129").Value = 2881.74315
$ws.Range("N129").Value = -12881.74315

$ws.Range("H137").Value = 1241.2609
$ws.Range("I137").Value = 920.5
$ws.Range("J137").Value = 2396
$ws.Range("K137").Value = 2761.5
$ws.Range("L137").Value = 7188
$ws.Range("M137").Value = -211.5
$ws.Range("N137").Value = -12288

$ws.Range("H138").Value = 1757.764
$ws.Range("I138").Value = 1171.75
$ws.Range("J138").Value = 1927.6232
$ws.Range("K138").Value = 3515.25
$ws.Range("L138").Value = 5782.8696
$ws.Range("M138").Value = 1624.75
$ws.Range("N138").Value = -16062.8696

$ws.Range("H141").Value = 730.8333
$ws.Range("I141").Value = 677
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 2031
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 3149
$ws.Range("N141").Value = -13360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -384
$ws.Range("N4").Value = -732

$ws.Range("H5").Value = 163.25
$ws.Range("I5").Value = 163.25
$ws.Range("K5").Value = 163.25
$ws.Range("M5").Value = -51.25

$ws.Range("H32").Value = 3732.7222
$ws.Range("I32").Value = 3967.6086
$ws.Range("J32").Value = 2382.125
$ws.Range("K32").Value = 3967.6086
$ws.Range("L32").Value = 2382.125
$ws.Range("M32").Value = -3680.6086
$ws.Range("N32").Value = -2956.125

$ws.Range("H61").Value = 1340.5
$ws.Range("I61").Value = 1340.5
$ws.Range("K61").Value = 1340.5
$ws.Range("M61").Value = -1128.5

$ws.Range("H97").Value = 350.63635
$ws.Range("I97").Value = 364.6
$ws.Range("J97").Value = 211
$ws.Range("K97").Value = 364.6
$ws.Range("L97").Value = 211
$ws.Range("M97").Value = 131.4
$ws.Range("N97").Value = -1203

$ws.Range("H122").Value = 1266.6666
$ws.Range("I122").Value = 1060
$ws.Range("K122").Value = 3180
$ws.Range("M122").Value = -730

$ws.Range("H132").Value = 2130.05
$ws.Range("I132").Value = 1741.3529
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 5224.0587
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -2694.0587
$ws.Range("N132").Value = -18057.9995

$ws.Range("H136").Value = 1340.5
$ws.Range("I136").Value = 1340.5
$ws.Range("K136").Value = 4021.5
$ws.Range("M136").Value = -1471.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 163.25
$ws.Range("I4").Value = 163.25
$ws.Range("K4").Value = 163.25
$ws.Range("M4").Value = -48.25

$ws.Range("H88").Value = 32000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 32000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 32000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -32812

$ws.Range("H91").Value = 32000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 32000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 32000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -34808

$ws.Range("H99").Value = 33334374
$ws.Range("I99").Value = 50000940
$ws.Range("K99").Value = 50000940
$ws.Range("M99").Value = -49999442

$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

$ws.Range("H134").Value = 7431.3687
$ws.Range("I134").Value = 1109
$ws.Range("K134").Value = 3327
$ws.Range("M134").Value = -792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1661.8182
$ws.Range("I58").Value = 1511.4286
$ws.Range("K58").Value = 1511.4286
$ws.Range("M58").Value = -1308.4286

$ws.Range("H105").Value = 799.875
$ws.Range("I105").Value = 733.3333
$ws.Range("J105").Value = 999.5
$ws.Range("K105").Value = 733.3333
$ws.Range("L105").Value = 999.5
$ws.Range("M105").Value = 1013.6667
$ws.Range("N105").Value = -4493.5

$ws.Range("H132").Value = 2420.7144
$ws.Range("I132").Value = 1198
$ws.Range("J132").Value = 2909.8
$ws.Range("K132").Value = 3594
$ws.Range("L132").Value = 8729.400000000001
$ws.Range("M132").Value = -1064
$ws.Range("N132").Value = -13789.4

$ws.Range("H136").Value = 1661.8182
$ws.Range("I136").Value = 1511.4286
$ws.Range("K136").Value = 4534.2858
$ws.Range("M136").Value = -1984.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4362.4287
$ws.Range("I80").Value = 1002
$ws.Range("J80").Value = 4922.5
$ws.Range("K80").Value = 3006
$ws.Range("L80").Value = 14767.5
$ws.Range("M80").Value = -2070
$ws.Range("N80").Value = -16639.5

$ws.Range("H83").Value = 4362.4287
$ws.Range("I83").Value = 1002
$ws.Range("J83").Value = 4922.5
$ws.Range("K83").Value = 9018
$ws.Range("L83").Value = 44302.5
$ws.Range("M83").Value = -4338
$ws.Range("N83").Value = -53662.5

$ws.Range("H86").Value = 405
$ws.Range("I86").Value = 300
$ws.Range("K86").Value = 900
$ws.Range("M86").Value = 286

$ws.Range("H89").Value = 405
$ws.Range("I89").Value = 300
$ws.Range("K89").Value = 2700
$ws.Range("M89").Value = 3228

$ws.Range("H114").Value = 452.2
$ws.Range("J114").Value = 496.5
$ws.Range("L114").Value = 1489.5
$ws.Range("N114").Value = -7997.5

$ws.Range("H130").Value = 2016.4
$ws.Range("J130").Value = 2016.4
$ws.Range("L130").Value = 6049.200000000001
$ws.Range("N130").Value = -16089.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 742.7143
$ws.Range("I97").Value = 742.7143
$ws.Range("K97").Value = 742.7143
$ws.Range("M97").Value = -246.7143

$ws.Range("H112").Value = 45046.1
$ws.Range("J112").Value = 45046.1
$ws.Range("L112").Value = 45046.1
$ws.Range("N112").Value = -47262.1

$ws.Range("H132").Value = 1792.1034
$ws.Range("I132").Value = 1461.0952
$ws.Range("J132").Value = 2661
$ws.Range("K132").Value = 4383.2856
$ws.Range("L132").Value = 7983
$ws.Range("M132").Value = -1853.2856
$ws.Range("N132").Value = -13043

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 689.3333
$ws.Range("I93").Value = 689.3333
$ws.Range("K93").Value = 689.3333
$ws.Range("M93").Value = 558.6667

$ws.Range("H98").Value = 8500
$ws.Range("J98").Value = 8500
$ws.Range("L98").Value = 8500
$ws.Range("N98").Value = -14490

$ws.Range("H100").Value = 2450
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2450
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2450
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3532

$ws.Range("H106").Value = 13999.5
$ws.Range("J106").Value = 13999.5
$ws.Range("L106").Value = 13999.5
$ws.Range("N106").Value = -16523.5

$ws.Range("H110").Value = 30644
$ws.Range("J110").Value = 30644
$ws.Range("L110").Value = 30644
$ws.Range("N110").Value = -38824

$ws.Range("H136").Value = 1233.0435
$ws.Range("I136").Value = 1159.762
$ws.Range("J136").Value = 2002.5
$ws.Range("K136").Value = 3479.286
$ws.Range("L136").Value = 6007.5
$ws.Range("M136").Value = -929.2860000000001
$ws.Range("N136").Value = -11107.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 16500
$ws.Range("J106").Value = 16500
$ws.Range("L106").Value = 16500
$ws.Range("N106").Value = -19024

$ws.Range("H122").Value = 12382429
$ws.Range("I122").Value = 16251369
$ws.Range("K122").Value = 48754107
$ws.Range("M122").Value = -48751657

$ws.Range("H132").Value = 2516.4707
$ws.Range("I132").Value = 2185.4666
$ws.Range("K132").Value = 6556.399800000001
$ws.Range("M132").Value = -4026.399800000001

$ws.Range("H136").Value = 666.2143
$ws.Range("I136").Value = 575.53845
$ws.Range("J136").Value = 1845
$ws.Range("K136").Value = 1726.61535
$ws.Range("L136").Value = 5535
$ws.Range("M136").Value = 823.38465
$ws.Range("N136").Value = -10635
